$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (shifts rows 13-23 down to 14-24)
$ws.Rows.Item(13).Insert()

# Row 13 gained a stray A13 cell from the insert (inherited style/content from row 12).
# Clear it out and reset its style to Normal so it matches an empty A13.
$ws.Cells.Item(13,1).ClearContents()
$ws.Cells.Item(13,1).Style = "Normal"

$ws.Range("B10").Value = "Apresentar os principais conceitos sobre as transformações de fases em materiais metálicos, poliméricos e cerâmicos abrangendo transformações difusionais e não-difusionais, a conceituação sobre nucleação e crescimento (aspectos energéticos) e sua relação com problemas práticos encontrados nas indústrias de processamento e de transformação de materiais."
$ws.Range("C10").Value = "Apresentar os principais conceitos sobre as transformações de fases em materiais metálicos, poliméricos e cerâmicos abrangendo transformações difusionais e não-difusionais, a conceituação sobre nucleação e crescimento (aspectos energéticos) e sua relação com problemas práticos encontrados nas indústrias de processamento e de transformação de materiais."
$ws.Range("B13").Value = "984972 - Hugo Ricardo Zschommler Sandim"
$ws.Range("C13").Value = "984972 - Hugo Ricardo Zschommler Sandim"
$ws.Range("B14").Value = "Difusão no estado sólido.Difusão em materiais não-metálicos (sólidos iônicos e polímeros).Recuperação, recristalização e crescimento de grão.Solidificação.Precipitação.Cinética de transformação no sistema Fe-C e em ligas não-ferrosas.Transformação de fases em vidros e cerâmicas.Transformação de fases em polímeros.Atividade experimental."
$ws.Range("C14").Value = "Difusão no estado sólido.Difusão em materiais não-metálicos (sólidos iônicos e polímeros).Recuperação, recristalização e crescimento de grão.Solidificação.Precipitação.Cinética de transformação no sistema Fe-C e em ligas não-ferrosas.Transformação de fases em vidros e cerâmicas.Transformação de fases em polímeros.Atividade experimental."
$ws.Range("B16").Value = "Introdução à difusão no estado sólido. Coeficiente de difusão. Leis de Fick. Difusão em soluções diluídas e na presença de um gradiente de concentração. Efeito Kirkendall.- Apresentar os fundamentos teóricos pertinentes à transformação de fases em materiais metálicos, cerâmicos e poliméricos.- Apresentar os conceitos fundamentais associados à nucleação, ao crescimento e à cinética de transformação de fases.- Descrições detalhadas de microestruturas fundidas e tratadas termicamente. Aspectos morfológicos relevantes.- Descrição das principais transformações de fase no estado sólido no sistema Fe-C e em algumas ligas não-ferrosas. Curvas TTT e CCT (TRC).- Estudar a transformação de fases durante o processamento termomecânico de metais e ligas.- Realização de prática experimental versando sobre tópicos da ementa."
$ws.Range("C16").Value = "Introdução à difusão no estado sólido. Coeficiente de difusão. Leis de Fick. Difusão em soluções diluídas e na presença de um gradiente de concentração. Efeito Kirkendall.- Apresentar os fundamentos teóricos pertinentes à transformação de fases em materiais metálicos, cerâmicos e poliméricos.- Apresentar os conceitos fundamentais associados à nucleação, ao crescimento e à cinética de transformação de fases.- Descrições detalhadas de microestruturas fundidas e tratadas termicamente. Aspectos morfológicos relevantes.- Descrição das principais transformações de fase no estado sólido no sistema Fe-C e em algumas ligas não-ferrosas. Curvas TTT e CCT (TRC).- Estudar a transformação de fases durante o processamento termomecânico de metais e ligas.- Realização de prática experimental versando sobre tópicos da ementa."
$ws.Range("B19").Value = "O aluno será avaliado ao longo do semestre por duas avaliações escritas (P1 e P2) e com pesos iguais."
$ws.Range("C19").Value = "O aluno será avaliado ao longo do semestre por duas avaliações escritas (P1 e P2) e com pesos iguais."
$ws.Range("B20").Value = "Nota Final NF = [P1 + P2]/2"
$ws.Range("C20").Value = "Nota Final NF = [P1 + P2]/2"
$ws.Range("B21").Value = "Para a recuperação será realizada uma prova escrita (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2"
$ws.Range("C21").Value = "Para a recuperação será realizada uma prova escrita (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2"
$ws.Range("B22").Value = "1 - Diffusion in solids. P.G. Shewmon, McGraw-Hill, 1963.2 - Phase transformation in metals. P.G. Shewmon, McGraw-Hill, 1969.3 - Recrystallization and related annealing phenomena. F.J. Humphreys and M. Hatherly, Pergamon, 1996.4 - Principles of solidification. B. Chalmers, Robert E. Krieger, 2nd. ed., 1977.5 - Precipitation hardening. A. Kelly, Pergamon, 1963.6 - Particle strengthening of metals and alloys. E. Nembach, John Wiley & Sons, 1997.7 - Propriedades dos materiais cerâmicos. L.H. Van Vlack, Edgard Blücher, 1973.8 - Textbook of polymer science. F.W. Billmeyer Jr., John Wiley & Sons, 1962.9 - Worked examples in the kinetics and thermodynamics of phase transformations. E.A. Wilson, The Institution of Metallurgists, s.d."
$ws.Range("C22").Value = "1 - Diffusion in solids. P.G. Shewmon, McGraw-Hill, 1963.2 - Phase transformation in metals. P.G. Shewmon, McGraw-Hill, 1969.3 - Recrystallization and related annealing phenomena. F.J. Humphreys and M. Hatherly, Pergamon, 1996.4 - Principles of solidification. B. Chalmers, Robert E. Krieger, 2nd. ed., 1977.5 - Precipitation hardening. A. Kelly, Pergamon, 1963.6 - Particle strengthening of metals and alloys. E. Nembach, John Wiley & Sons, 1997.7 - Propriedades dos materiais cerâmicos. L.H. Van Vlack, Edgard Blücher, 1973.8 - Textbook of polymer science. F.W. Billmeyer Jr., John Wiley & Sons, 1962.9 - Worked examples in the kinetics and thermodynamics of phase transformations. E.A. Wilson, The Institution of Metallurgists, s.d."
